$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.538.43"
$ws.Range("E2").Value = "  +2.05%  "
$ws.Range("D3").Value = "1.672.41"
$ws.Range("E3").Value = "  +2.40%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.64"
$ws.Range("E5").Value = "  +2.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.528"
$ws.Range("E6").Value = "  +2.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "29.59"
$ws.Range("E8").Value = "  +3.62%  "
$ws.Range("E9").Value = "  +2.78%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0643"
$ws.Range("E10").Value = "  +5.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0904"
$ws.Range("E11").Value = "  -0.51%  "
$ws.Range("D12").Value = "1.912.06"
$ws.Range("E12").Value = "  +2.38%  "
$ws.Range("E13").Value = "  +9.03%  "
$ws.Range("D14").Value = "1.669.37"
$ws.Range("E14").Value = "  +2.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "10.19"
$ws.Range("E15").Value = "  +9.78%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.00"
$ws.Range("E16").Value = "  +4.01%  "
$ws.Range("D17").Value = "30.560.30"
$ws.Range("E17").Value = "  +2.10%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "66.37"
$ws.Range("E18").Value = "  +3.63%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.70"
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("D20").Value = "0.0₃0724"
$ws.Range("E20").Value = "  +3.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.999"
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("E22").Value = "  +3.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.00"
$ws.Range("E23").Value = "  +2.18%  "
$ws.Range("E24").Value = "  -0.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.47"
$ws.Range("E25").Value = "  +0.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.86"
$ws.Range("E26").Value = "  +2.21%  "
$ws.Range("E27").Value = "  +2.58%  "
$ws.Range("E28").Value = "  +1.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("E30").Value = "  +1.81%  "
$ws.Range("E31").Value = "  +2.77%  "
$ws.Range("E32").Value = "  +2.77%  "
$ws.Range("E33").Value = "  +3.28%  "
$ws.Range("D34").Value = "1.493.81"
$ws.Range("E34").Value = "  +4.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.76"
$ws.Range("E35").Value = "  +7.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "84.79"
$ws.Range("E36").Value = "  +12.46%  "
$ws.Range("E37").Value = "  -0.64%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.598"
$ws.Range("E38").Value = "  +8.16%  "
$ws.Range("E39").Value = "  +5.17%  "
$ws.Range("E40").Value = "  -4.59%  "
$ws.Range("E41").Value = "  -0.24%  "
$ws.Range("E42").Value = "  +1.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0498"
$ws.Range("E43").Value = "  +1.37%  "
$ws.Range("E44").Value = "  -1.60%  "
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.998"
$ws.Range("E46").Value = "  -0.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.51"
$ws.Range("E47").Value = "  +2.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "50.82"
$ws.Range("E48").Value = "  -3.57%  "
$ws.Range("D49").Value = "1.805.38"
$ws.Range("E49").Value = "  +1.73%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "94.84"
$ws.Range("E50").Value = "  +4.98%  "
$ws.Range("E51").Value = "  -1.16%  "
